$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the changed values in row 2 (new environment, policy number, incident date) ---
$ws.Range("E2").Value = 12112002070
$ws.Range("G2").Value = "'23/06/2022"
$ws.Range("A2").Value = "ssurgwsoadev4-oci.opc.oracleoutsourcing.com"
$ws.Range("B2").Value = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/cc/ClaimCenter.do"

# Keep a copy of the existing hyperlink cell format (font/underline) so that
# rebuilding the hyperlinks below does not invent a brand-new cell style.
$ws.Range("B2").Copy()

# The engine's Hyperlinks.Delete() removes the whole sheet's hyperlinks no
# matter which range it is invoked from, so rebuild the full collection -
# this mirrors what Excel itself does internally when a single hyperlink
# (B2, pointing at the old "ssurgwsoadev4" host) is edited: every other
# hyperlink relationship gets renumbered too.
$ws.Range("A1").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("J2"), "mailto:aseguradosgw@gmail.com")
$ws.Range("J2").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("B3"), "https://i-preproducciongestion.segurossura.com.ar/cc/ClaimCenter.do")
$ws.Range("B3").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("J3"), "mailto:aseguradosgw@gmail.com")
$ws.Range("J3").PasteSpecial(-4122)

$ws.Hyperlinks.Add($ws.Range("B2"), "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/cc/ClaimCenter.do")
$ws.Range("B2").PasteSpecial(-4122)

# --- Move the active selection to E3, matching the saved sheet view ---
$ws.Range("E3").Select()
